$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the week's day-labels (column A, rows 11-17) to the new week (w/c Sun 04/05)
$ws.Range("A11").Value = "Sun 04/05"
$ws.Range("A12").Value = "Mon 05/05"
$ws.Range("A13").Value = "Tue 06/05"
$ws.Range("A14").Value = "Wed  07/05"
$ws.Range("A15").Value = "Thur  08/05"
$ws.Range("A16").Value = "Fri   09/05"
$ws.Range("A17").Value = "Sat  10/05"

# Update "Week of:" date (G8) to the Sunday of that week
$ws.Range("G8").Value = 41763

# Update the selected range / view to match where the user left off editing
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("G8:H8").Select()
